$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44391
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 21000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21500
$ws.Range("P2").Value = 538

# Row 3
$ws.Range("D3").Value = 44489
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("P3").Value = 338

# Row 4
$ws.Range("D4").Value = 44426
$ws.Range("J4").Value = 150

# Row 5
$ws.Range("D5").Value = 44405
$ws.Range("K5").Value = 21000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 21500
$ws.Range("P5").Value = 538

# Row 6
$ws.Range("D6").Value = 44370
$ws.Range("H6").Value = 'Argentina(o)'
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 21000
$ws.Range("M6").Value = 20429
$ws.Range("N6").Value = '$/caja 50 unidades'
$ws.Range("P6").Value = 409
$ws.Range("Q6").Value = 50

# Row 7
$ws.Range("D7").Value = 44370
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 23000
$ws.Range("M7").Value = 22500
$ws.Range("P7").Value = 562

# Row 8
$ws.Range("D8").Value = 44419
$ws.Range("H8").Value = 'Symphony'
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 21000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21500
$ws.Range("P8").Value = 430

# Row 9
$ws.Range("D9").Value = 44706
$ws.Range("H9").Value = 'Madrigal'
$ws.Range("J9").Value = 250
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 538
$ws.Range("Q9").Value = 40

# Row 10
$ws.Range("D10").Value = 44482
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("P10").Value = 362

# Row 11
$ws.Range("D11").Value = 44483
$ws.Range("H11").Value = 'Madrigal'
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("N11").Value = '$/caja 40 unidades'
$ws.Range("P11").Value = 362
$ws.Range("Q11").Value = 40

# Row 12
$ws.Range("D12").Value = 44167
$ws.Range("H12").Value = 'Española'
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 13500
$ws.Range("N12").Value = '$/caja 30 unidades'
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 450
$ws.Range("Q12").Value = 30

# Row 13
$ws.Range("D13").Value = 44384
$ws.Range("H13").Value = 'Madrigal'
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 21000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 21500
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("P13").Value = 538
$ws.Range("Q13").Value = 40

# Row 14
$ws.Range("D14").Value = 44384
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 19000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19333
$ws.Range("N14").Value = '$/caja 50 unidades'
$ws.Range("P14").Value = 387
$ws.Range("Q14").Value = 50

# Row 15
$ws.Range("H15").Value = 'Symphony'
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 21000
$ws.Range("M15").Value = 20400
$ws.Range("P15").Value = 510

# Row 16
$ws.Range("D16").Value = 44468
$ws.Range("H16").Value = 'Argentina(o)'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17500
$ws.Range("P16").Value = 350

# Row 17
$ws.Range("D17").Value = 44356
$ws.Range("H17").Value = 'Argentina(o)'
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 19000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19500
$ws.Range("N17").Value = '$/caja 50 unidades'
$ws.Range("P17").Value = 390
$ws.Range("Q17").Value = 50

# Row 18
$ws.Range("D18").Value = 44412
$ws.Range("H18").Value = 'Symphony'
$ws.Range("J18").Value = 240

# Row 19
$ws.Range("D19").Value = 44435
$ws.Range("H19").Value = 'Madrigal'
$ws.Range("K19").Value = 19000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19500
$ws.Range("N19").Value = '$/caja 40 unidades'
$ws.Range("O19").Value = 'Región de Coquimbo'
$ws.Range("P19").Value = 488
$ws.Range("Q19").Value = 40

# Row 20
$ws.Range("D20").Value = 44433
$ws.Range("H20").Value = 'Madrigal'
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 19000
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 19500
$ws.Range("P20").Value = 488

# Row 21
$ws.Range("D21").Value = 44398
$ws.Range("J21").Value = 170
$ws.Range("K21").Value = 21000
$ws.Range("L21").Value = 22000
$ws.Range("M21").Value = 21500
$ws.Range("P21").Value = 538

# Row 23
$ws.Range("D23").Value = 44363
$ws.Range("J23").Value = 160
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19500
$ws.Range("P23").Value = 488

# Row 24
$ws.Range("H24").Value = 'Madrigal'
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 21000
$ws.Range("M24").Value = 20333
$ws.Range("P24").Value = 508

# Row 25
$ws.Range("D25").Value = 44377
$ws.Range("H25").Value = 'Symphony'
$ws.Range("J25").Value = 60
